$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") is updated from 2026-02-24 to 2026-02-25 for every data row (2-9).
$ws.Range("C2:C9").Value = 46078

# Rows 3-9 are re-ordered; only columns A (Beteckning), B (Datum) and G (Area (ha))
# differ between rows, so rewrite those three columns for each row to match the
# new order.
$rows = @(
    @{ Row = 3; A = "A 34310-2024"; B = 45524; G = 4.8 },
    @{ Row = 4; A = "A 25617-2024"; B = 45463; G = 2.3 },
    @{ Row = 5; A = "A 45983-2023"; B = 45196; G = 0.6 },
    @{ Row = 6; A = "A 54782-2022"; B = 44883; G = 5.5 },
    @{ Row = 7; A = "A 843-2024";   B = 45300; G = 0.8 },
    @{ Row = 8; A = "A 844-2024";   B = 45300; G = 1.2 },
    @{ Row = 9; A = "A 17908-2021"; B = 44301; G = 0.9 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
